# Milestone M0 doc update: server moved from 50.18.130.129 to 13.57.40.53,
# and the Mongo Atlas username changed from "csc@890" to "csc890team4".

$d = $word.ActiveDocument

# 1) Website URL table row: http://50.18.130.129:3000/ -> http://13.57.40.53:3000/
$d.Content.Find.Execute(
    "http://50.18.130.129:3000/", $true, $false, $false, $false, $false,
    $true, 1, $false, "http://13.57.40.53:3000/", 2)

# 2) SSH URL table row: ec2-50-18-130-129 -> ec2-13-57-40-53
$d.Content.Find.Execute(
    "ubuntu@ec2-50-18-130-129.us-west-1.compute.amazonaws.com", $true, $false, $false, $false, $false,
    $true, 1, $false, "ubuntu@ec2-13-57-40-53.us-west-1.compute.amazonaws.com", 2)

# 3) Mongo connection string table row: //50.18.130.129: -> //13.57.40.53:
$d.Content.Find.Execute(
    "mongodb://50.18.130.129:27017", $true, $false, $false, $false, $false,
    $true, 1, $false, "mongodb://13.57.40.53:27017", 2)

# 4) Mongo database password table row: csc@890 -> csc890team4
$d.Content.Find.Execute(
    "csc@890", $true, $false, $false, $false, $false,
    $true, 1, $false, "csc890team4", 2)
